# Apply the Feb 26 2024 cryptos-list refresh (prices / 1h volume, plus the
# three coin-row swaps and the Mantle -> ThetaToken replacement at row 51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values that render as plain numbers (e.g. "0.999", "1.00", "38.20") must be
# written with a leading apostrophe so Excel stores/displays them as literal
# text instead of silently converting them to numbers and dropping trailing
# zeros / the decimal formatting used by the source data feed.

$ws.Range("D2").Value = '52.740.62'
$ws.Range("E2").Value = '  +2.38%  '

$ws.Range("D3").Value = '3.125.64'
$ws.Range("E3").Value = '  +2.70%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''395.63'
$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("D6").Value = '''104.54'
$ws.Range("E6").Value = '  +1.82%  '

$ws.Range("D7").Value = '''0.541'
$ws.Range("E7").Value = '  -0.41%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").Value = '''0.606'
$ws.Range("E9").Value = '  +3.95%  '

$ws.Range("D10").Value = '''38.20'
$ws.Range("E10").Value = '  +3.66%  '

$ws.Range("E11").Value = '  +0.97%  '

$ws.Range("D12").Value = '''0.0865'
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").Value = '3.607.93'
$ws.Range("E13").Value = '  +2.13%  '

$ws.Range("D14").Value = '''18.83'
$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").Value = '''7.85'
$ws.Range("E15").Value = '  +1.64%  '

$ws.Range("E16").Value = '  +8.16%  '

$ws.Range("D17").Value = '3.121.84'
$ws.Range("E17").Value = '  +2.38%  '

$ws.Range("D18").Value = '''10.76'
$ws.Range("E18").Value = '  +2.54%  '

$ws.Range("D19").Value = '52.562.05'
$ws.Range("E19").Value = '  +1.98%  '

$ws.Range("D20").Value = '''3.24'
$ws.Range("E20").Value = '  +3.78%  '

$ws.Range("E21").Value = '  +3.02%  '

$ws.Range("D22").Value = '0.0₃0972'
$ws.Range("E22").Value = '  +1.05%  '

$ws.Range("D23").Value = '''71.09'
$ws.Range("E23").Value = '  +1.34%  '

$ws.Range("D24").Value = '''269.36'
$ws.Range("E24").Value = '  +0.54%  '

$ws.Range("D25").Value = '''3.22'
$ws.Range("E25").Value = '  +2.34%  '

$ws.Range("D26").Value = '''8.06'
$ws.Range("E26").Value = '  -3.35%  '

$ws.Range("D27").Value = '''27.59'
$ws.Range("E27").Value = '  +2.10%  '

$ws.Range("D28").Value = '''7.48'
$ws.Range("E28").Value = '  +3.60%  '

$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  -2.20%  '

$ws.Range("E31").Value = '  +0.16%  '

$ws.Range("D32").Value = '''10.90'
$ws.Range("E32").Value = '  +6.17%  '

$ws.Range("D33").Value = '''36.90'
$ws.Range("E33").Value = '  +7.11%  '

$ws.Range("E34").Value = '  +10.28%  '

$ws.Range("E35").Value = '  +0.74%  '

$ws.Range("D36").Value = '''50.10'
$ws.Range("E36").Value = '  -0.77%  '

$ws.Range("D37").Value = '''0.999'
$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").Value = '''3.43'
$ws.Range("E38").Value = '  +2.07%  '

$ws.Range("D39").Value = '''4.11'
$ws.Range("E39").Value = '  +10.75%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''2.69'
$ws.Range("E40").Value = '  +6.55%  '

$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '''0.293'
$ws.Range("E41").Value = '  +1.61%  '

$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").Value = '''17.00'
$ws.Range("E42").Value = '  +0.28%  '

$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").Value = '''130.34'
$ws.Range("E43").Value = '  +1.78%  '

$ws.Range("D44").Value = '''1.88'
$ws.Range("E44").Value = '  +0.34%  '

$ws.Range("E45").Value = '  +0.48%  '

$ws.Range("D46").Value = '''22.23'
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("E47").Value = '  -1.35%  '

$ws.Range("E48").Value = '  -0.55%  '

$ws.Range("D49").Value = '2.082.85'
$ws.Range("E49").Value = '  +2.22%  '

$ws.Range("D50").Value = '''0.0527'
$ws.Range("E50").Value = '  +33.67%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '''1.71'
$ws.Range("E51").Value = '  +28.00%  '
